$wb = $excel.ActiveWorkbook

# Rename the existing sheet to "Towers"
$towers = $wb.Worksheets.Item(1)
$towers.Name = "Towers"

# Add a new sheet after "Towers" named "Feuil1"
$feuil1 = $wb.Worksheets.Add($null, $towers)
$feuil1.Name = "Feuil1"

# --- Towers sheet content ---
$towers.Range("B1").Value = "Goblin"
$towers.Range("C1").Value = "Eye"
$towers.Range("D1").Value = "Mushroom"

$towers.Range("A2").Value = "Damage"
$towers.Range("B2").Value = 10
$towers.Range("C2").Value = 75

$towers.Range("A3").Value = "Cost"
$towers.Range("B3").Value = 50
$towers.Range("C3").Value = 50
$towers.Range("D3").Value = 50

$towers.Range("A4").Value = "UpgradeCost"
$towers.Range("B4").Value = "50/70"
$towers.Range("C4").Value = "50/70"
$towers.Range("D4").Value = "50/70"

$towers.Range("A5").Value = "Experience"
$towers.Range("B5").Value = 2
$towers.Range("C5").Value = 2
$towers.Range("D5").Value = 2

$towers.Range("A8").Value = "NbrSprite"
$towers.Range("B8").Value = 15
$towers.Range("C8").Value = 26

$towers.Range("A9").Value = "Vitesse"
$towers.Range("B9").Value = 15
$towers.Range("C9").Value = 12

$towers.Range("A10").Value = "dps"
$towers.Range("B10").Formula = "=B2*(B9/B8)"
$towers.Range("C10").Formula = "=C2*(C9/C8)"
$towers.Range("C10").NumberFormat = "0"

$towers.Columns.Item(1).ColumnWidth = 16.5
$towers.Columns.Item(3).ColumnWidth = 12.64

# --- Feuil1 sheet content ---
$feuil1.Range("A1").Value = "Vitesse"
$feuil1.Range("A2").Value = "HP"
$feuil1.Range("A3").Value = "Experience"

$null = $towers.Range("F11").Select()
$null = $feuil1.Activate()
$null = $feuil1.Range("A4").Select()
